$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.346.93"
$ws.Range("E2").Value = "  +1.27%  "
$ws.Range("D3").Value = "3.137.05"
$ws.Range("E3").Value = "  +1.59%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'604.04"
$ws.Range("E5").Value = "  -0.41%  "
$ws.Range("D6").Value = "'143.03"
$ws.Range("E6").Value = "  -0.69%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "3.134.12"
$ws.Range("E8").Value = "  +1.62%  "
$ws.Range("D9").Value = "'0.523"
$ws.Range("E9").Value = "  +1.21%  "
$ws.Range("E10").Value = "  +1.51%  "
$ws.Range("D11").Value = "'5.38"
$ws.Range("E11").Value = "  +4.53%  "
$ws.Range("E12").Value = "  +0.45%  "
$ws.Range("D13").Value = "'0.0000256"
$ws.Range("E13").Value = "  +4.34%  "
$ws.Range("D14").Value = "'35.23"
$ws.Range("E14").Value = "  +1.04%  "
$ws.Range("D15").Value = "3.658.15"
$ws.Range("E15").Value = "  +1.68%  "
$ws.Range("E16").Value = "  +3.19%  "
$ws.Range("D17").Value = "64.174.49"
$ws.Range("E17").Value = "  +0.90%  "
$ws.Range("D18").Value = "3.154.64"
$ws.Range("E18").Value = "  +2.18%  "
$ws.Range("D19").Value = "'6.86"
$ws.Range("E19").Value = "  +1.71%  "
$ws.Range("D20").Value = "'478.36"
$ws.Range("E20").Value = "  +1.53%  "
$ws.Range("D21").Value = "'14.57"
$ws.Range("E21").Value = "  +0.74%  "
$ws.Range("E22").Value = "  +2.02%  "
$ws.Range("D23").Value = "'7.68"
$ws.Range("E23").Value = "  +0.78%  "
$ws.Range("D24").Value = "'85.10"
$ws.Range("E24").Value = "  +2.57%  "
$ws.Range("D25").Value = "'13.40"
$ws.Range("E25").Value = "  -0.09%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("E27").Value = "  -0.25%  "
$ws.Range("D28").Value = "'8.39"
$ws.Range("E28").Value = "  +1.82%  "
$ws.Range("D29").Value = "'7.21"
$ws.Range("E29").Value = "  +9.02%  "
$ws.Range("E30").Value = "  -3.83%  "
$ws.Range("E31").Value = "  +2.87%  "
$ws.Range("B32").Value = "FirstDigitalUSD"
$ws.Range("C32").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D32").Value = "'1.00"
$ws.Range("E32").Value = "  +0.04%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").Value = "'26.92"
$ws.Range("E33").Value = "  +3.88%  "
$ws.Range("D34").Value = "'2.64"
$ws.Range("E34").Value = "  -2.22%  "
$ws.Range("E35").Value = "  +0.77%  "
$ws.Range("D36").Value = "0.0₃0772"
$ws.Range("E36").Value = "  +7.63%  "
$ws.Range("D37").Value = "'5.96"
$ws.Range("E37").Value = "  +1.58%  "
$ws.Range("D38").Value = "'52.43"
$ws.Range("E38").Value = "  +0.83%  "
$ws.Range("D39").Value = "'3.02"
$ws.Range("E39").Value = "  +5.46%  "
$ws.Range("D40").Value = "'444.72"
$ws.Range("E40").Value = "  -2.19%  "
$ws.Range("E41").Value = "  +1.01%  "
$ws.Range("E42").Value = "  +0.86%  "
$ws.Range("D43").Value = "'8.20"
$ws.Range("E43").Value = "  -0.65%  "
$ws.Range("D44").Value = "2.852.41"
$ws.Range("E44").Value = "  +1.60%  "
$ws.Range("E45").Value = "  -0.44%  "
$ws.Range("D46").Value = "'2.22"
$ws.Range("E46").Value = "  +0.56%  "
$ws.Range("E47").Value = "  +2.09%  "
$ws.Range("E48").Value = "  +0.04%  "
$ws.Range("D49").Value = "'25.99"
$ws.Range("E49").Value = "  +0.93%  "
$ws.Range("E50").Value = "  +0.73%  "
$ws.Range("D51").Value = "'119.78"
$ws.Range("E51").Value = "  +1.74%  "
